# Set H26:H176 (Industries column) from 1 to 0, matching the commit's
# update of policy data rows (one cell per day, rows 26-176).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H26:H176").Value = 0
